$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Power comparisons for 3x2 - update computed values in the LP1/LP2 comparison table

$ws.Range("J2:L2").Value = "0.0077 -0.0802"
$ws.Range("G3").Value = "0.535 5e-04"
$ws.Range("J3:L3").Value = "0.0291 -0.0475"

$ws.Range("G4:G6").Value = "0.9802 0.0485"

$ws.Range("J4:L4").Value = "0.2259 5e-04"
$ws.Range("J5:L5").Value = "0.2259 5e-04"
$ws.Range("J6:L6").Value = "0.2259 5e-04"
$ws.Range("J8:L8").Value = "0.2259 5e-04"
$ws.Range("J9:L9").Value = "0.2259 5e-04"

$ws.Range("M4:O4").Value = "0.3685 0"
$ws.Range("M5:O5").Value = "0.3685 0"
$ws.Range("M6:O6").Value = "0.3685 0"
$ws.Range("M8:O8").Value = "0.3685 0"
$ws.Range("M9:O9").Value = "0.3685 0"

$ws.Range("H7:I7").Value = "< -0.0485"
$ws.Range("M7:O7").Value = "< -0.0485"
$ws.Range("J7:L7").Value = "0.0188 -0.048"

$ws.Range("M10:O10").Value = "0.655 -5e-04"
$ws.Range("M11:O11").Value = "0.655 -5e-04"
$ws.Range("M12:O12").Value = "0.655 -5e-04"
